# The sheet's table (header "Enero..Total mes" in row 1, 15 "Vendedor"
# rows, and a "Total mes" footer row in row 17) is pushed down by two rows:
# everything that used to live in rows 1-17 now lives in rows 3-19, and the
# sheet's used range/dimension becomes A3:H19 (was A1:H17). Inserting two
# whole rows above row 1 reproduces exactly that shift, carrying every
# cell's value/style/number-format down with it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:2").Insert()

# The saved selection also moves from C24 to the "Total mes" column for the
# whole data range (H4:H18), matching the post-edit sheetView.
[void]$ws.Range("H4:H18").Select()
